$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1
#    title paragraph ("Play Bubble Craze for Free: A Unique Slot Game by
#    IGT"). We use Range.InsertXML so the new paragraph is created without
#    inheriting the Heading1 paragraph style (it ends up as a plain body
#    paragraph, matching the rest of the document).
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End)

$openXml = $d.Content.WordOpenXML
$bodyTagEnd = $openXml.IndexOf("<w:body>") + 8
$docTagStart = $openXml.IndexOf("<w:document ")
$docHeader = $openXml.Substring($docTagStart, $bodyTagEnd - $docTagStart)

$packageHead = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>'
$packageTail = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$replacementParas = '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Bubble Craze for Free: A Unique Slot Game by IGT</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Try Bubble Craze for free and experience a unique approach to slot machine gameplay with excellent graphics and bonus features by trusted developer IGT.</w:t></w:r></w:p>'

$flatOpcXml = $packageHead + $docHeader + $replacementParas + $packageTail
$titleRange.InsertXML($flatOpcXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Near the end of the document, remove the duplicated bold title
#    paragraph ("Play Bubble Craze for Free: A Unique Slot Game by IGT")
#    that used to sit before the italic meta-description paragraph. (The
#    very first paragraph has the same text but a Heading1 style, so only
#    delete the later, non-heading occurrence. Paragraph.Range.Text always
#    carries a trailing paragraph-mark character, so trim before comparing.)
# ---------------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
for ($i = $paraCount; $i -ge 2; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq "Play Bubble Craze for Free: A Unique Slot Game by IGT") {
        $p.Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# 3) Replace the italic tagline/meta text at the very end of the document
#    with the new image-prompt text. The search is scoped to just the last
#    paragraph so it cannot also match the identical sentence that now
#    lives inside the new "Meta description" paragraph near the top.
# ---------------------------------------------------------------------------
$oldText = "Try Bubble Craze for free and experience a unique approach to slot machine gameplay with excellent graphics and bonus features by trusted developer IGT."
$newText = "Create a feature image for Bubble Craze that showcases a happy Maya warrior with glasses in a cartoon style. The background of the image should resemble a blurred image of bubbles, and the warrior should be holding a bubble wand. The warrior should be surrounded by bubbles of different colors, and there should be a transformation and multiplier bubble present in the image."

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

Write-Output "Edit complete"
